# Update cryptos list data to reflect refreshed prices/volumes (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every cell in columns B-E is stored as text (inlineStr) in this workbook, including
# numeric-looking Price values like "227.34". A leading apostrophe (quote prefix) is
# prepended so Excel keeps each new value as literal text instead of converting it to
# a number/date.
$ws.Range("D2").Value = "'" + '37.674.79'
$ws.Range("E2").Value = "'" + '  -1.63%  '
$ws.Range("D3").Value = "'" + '2.025.27'
$ws.Range("E3").Value = "'" + '  -2.12%  '
$ws.Range("E4").Value = "'" + '  -0.02%  '
$ws.Range("D5").Value = "'" + '227.34'
$ws.Range("E5").Value = "'" + '  -1.58%  '
$ws.Range("D6").Value = "'" + '0.602'
$ws.Range("E6").Value = "'" + '  -2.92%  '
$ws.Range("D7").Value = "'" + '59.66'
$ws.Range("E7").Value = "'" + '  -3.48%  '
$ws.Range("E8").Value = "'" + '  +0.06%  '
$ws.Range("D9").Value = "'" + '0.374'
$ws.Range("E9").Value = "'" + '  -3.89%  '
$ws.Range("E10").Value = "'" + '  +0.86%  '
$ws.Range("E11").Value = "'" + '  -0.82%  '
$ws.Range("D12").Value = "'" + '2.326.39'
$ws.Range("E12").Value = "'" + '  -2.04%  '
$ws.Range("D13").Value = "'" + '14.36'
$ws.Range("E13").Value = "'" + '  -4.11%  '
$ws.Range("E14").Value = "'" + '  -2.50%  '
$ws.Range("D15").Value = "'" + '0.763'
$ws.Range("E15").Value = "'" + '  -0.74%  '
$ws.Range("E16").Value = "'" + '  -3.21%  '
$ws.Range("D17").Value = "'" + '2.031.37'
$ws.Range("E17").Value = "'" + '  -1.80%  '
$ws.Range("D18").Value = "'" + '37.655.92'
$ws.Range("E18").Value = "'" + '  -1.51%  '
$ws.Range("E19").Value = "'" + '  -1.46%  '
$ws.Range("D20").Value = "'" + '5.87'
$ws.Range("E20").Value = "'" + '  -7.04%  '
$ws.Range("E21").Value = "'" + '  -2.15%  '
$ws.Range("D22").Value = "'" + '223.54'
$ws.Range("E22").Value = "'" + '  -1.55%  '
$ws.Range("E23").Value = "'" + '  +0.10%  '
$ws.Range("E24").Value = "'" + '  -0.72%  '
$ws.Range("D25").Value = "'" + '2.24'
$ws.Range("E25").Value = "'" + '  -0.81%  '
$ws.Range("D26").Value = "'" + '167.65'
$ws.Range("E26").Value = "'" + '  +0.54%  '
$ws.Range("D27").Value = "'" + '9.29'
$ws.Range("E27").Value = "'" + '  -0.53%  '
$ws.Range("E28").Value = "'" + '  -3.41%  '
$ws.Range("D29").Value = "'" + '18.76'
$ws.Range("E29").Value = "'" + '  -2.26%  '
$ws.Range("D30").Value = "'" + '1.25'
$ws.Range("E30").Value = "'" + '  -6.66%  '
$ws.Range("E31").Value = "'" + '  -0.53%  '
$ws.Range("D32").Value = "'" + '2.21'
$ws.Range("E32").Value = "'" + '  +7.67%  '
$ws.Range("D33").Value = "'" + '4.37'
$ws.Range("E33").Value = "'" + '  -4.76%  '
$ws.Range("D34").Value = "'" + '0.0602'
$ws.Range("E34").Value = "'" + '  -1.13%  '
$ws.Range("E35").Value = "'" + '  -4.23%  '
$ws.Range("D36").Value = "'" + '6.40'
$ws.Range("E36").Value = "'" + '  +1.83%  '
$ws.Range("D37").Value = "'" + '2.29'
$ws.Range("E37").Value = "'" + '  -2.14%  '
$ws.Range("D38").Value = "'" + '3.38'
$ws.Range("E38").Value = "'" + '  +1.51%  '
$ws.Range("E39").Value = "'" + '  -0.09%  '
$ws.Range("D40").Value = "'" + '17.82'
$ws.Range("E40").Value = "'" + '  +3.12%  '
$ws.Range("D41").Value = "'" + '1.529.98'
$ws.Range("E41").Value = "'" + '  -0.02%  '
$ws.Range("E42").Value = "'" + '  -1.87%  '
$ws.Range("D43").Value = "'" + '95.33'
$ws.Range("E43").Value = "'" + '  -3.22%  '
$ws.Range("E44").Value = "'" + '  -2.21%  '
$ws.Range("D45").Value = "'" + '0.0906'
$ws.Range("E45").Value = "'" + '  -2.95%  '
$ws.Range("B46").Value = "'" + 'TrustWalletToken'
$ws.Range("C46").Value = "'" + 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").Value = "'" + '1.10'
$ws.Range("E46").Value = "'" + '  -3.06%  '
$ws.Range("B47").Value = "'" + 'FTXToken'
$ws.Range("C47").Value = "'" + 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = "'" + '4.05'
$ws.Range("E47").Value = "'" + '  +0.27%  '
$ws.Range("B48").Value = "'" + 'MXToken'
$ws.Range("C48").Value = "'" + 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D48").Value = "'" + '2.96'
$ws.Range("E48").Value = "'" + '  +0.12%  '
$ws.Range("B49").Value = "'" + 'ARBITRUM'
$ws.Range("C49").Value = "'" + 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = "'" + '1.00'
$ws.Range("E49").Value = "'" + '  -2.66%  '
$ws.Range("E50").Value = "'" + '  -0.30%  '
$ws.Range("D51").Value = "'" + '2.216.42'
$ws.Range("E51").Value = "'" + '  -2.01%  '
